$wb = $excel.ActiveWorkbook

# --- Insert the new "Floor Features" worksheet right after "Floors" ---
# Worksheet references returned by Add() track POSITION, not identity, so
# adding a second sheet at the same anchor point shifts any earlier
# reference. To land the new sheet with the sheetId the target file expects
# (9, i.e. one past the highest sheetId ever used, 8), we first add a
# disposable placeholder sheet at the very end of the workbook (consuming
# sheetId 8) and only then add the real "Floor Features" sheet after
# "Floors" (consuming sheetId 9). The placeholder is deleted afterwards.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wb.Worksheets.Add($null, $lastSheet, 1) | Out-Null

$floorsSheet = $wb.Worksheets.Item("Floors")
$newSheet = $wb.Worksheets.Add($null, $floorsSheet, 1)
$newSheet.Name = "Floor Features"

# The placeholder always ends up as the very last sheet (it was inserted
# after the then-last sheet, and the subsequent insert happened earlier in
# the tab order, which doesn't touch tab positions after it). Re-resolve it
# fresh by position right before deleting instead of trusting any earlier
# handle, since handles track position and go stale across further Add()s.
$wb.Worksheets.Item($wb.Worksheets.Count).Delete() | Out-Null

# --- Populate "Floor Features" ---
# Write order matters for shared-string allocation order, so values are
# poked in the same sequence the source workbook used.
$ff = $wb.Worksheets.Item("Floor Features")

$ff.Range("A1").Value = "Elevation Title"
$ff.Range("B1").Value = "Floor Title"
$ff.Range("D1").Value = "Feature Title"
$ff.Range("E1").Value = "Feature Price"

$ff.Range("A2").Value = "Elante"
$ff.Range("B2").Value = "Basement"
$ff.Range("C2").Value = 1

$ff.Range("A3").Value = "Elante"
$ff.Range("B3").Value = "Basement"
$ff.Range("C3").Value = 2
$ff.Range("D3").Value = "gourment(100)"
$ff.Range("E3").Value = 4000

$ff.Range("A4").Value = "Elante"
$ff.Range("B4").Value = "Basement"
$ff.Range("C4").Value = 2
$ff.Range("D4").Value = "gourment(10007)"
$ff.Range("E4").Value = 1000

$ff.Range("D2").Value = "KITCHEN"

$ff.Range("C1").Value = "(feature_group=1, feature=2)"

$ff.Range("A5").Value = "Elante"
$ff.Range("B5").Value = "Basement"
$ff.Range("C5").Value = 1
$ff.Range("D5").Value = "GARAGE"

$ff.Range("A1:E1").Font.Bold = $true

# --- Selection tweaks on existing sheets ---
$elevations = $wb.Worksheets.Item("Elevations")
$elevations.Range("A2").Select() | Out-Null

$floors = $wb.Worksheets.Item("Floors")
$floors.Range("A2").Select() | Out-Null

# --- Final active sheet/selection: "Floor Features", cell D5 ---
$ff.Range("D5").Select() | Out-Null
